$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks like a number but must stay text
# (preserves formatting such as trailing zeros, e.g. "1.00").
$textCells = @("D24", "D31", "D45")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.481.52"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "2.362.24"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "312.82"
$ws.Range("E5").Value = "  +5.49%  "
$ws.Range("D6").Value = "109.25"
$ws.Range("D7").Value = "0.641"
$ws.Range("E7").Value = "  +3.56%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.643"
$ws.Range("E9").Value = "  +6.89%  "
$ws.Range("D10").Value = "43.03"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "8.83"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Value = "16.51"
$ws.Range("E15").Value = "  +10.02%  "
$ws.Range("D16").Value = "2.716.37"
$ws.Range("E16").Value = "  +6.26%  "
$ws.Range("D17").Value = "2.419.81"
$ws.Range("E17").Value = "  +8.60%  "
$ws.Range("D18").Value = "43.465.68"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").Value = "7.24"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").Value = "75.23"
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("D23").Value = "2.55"
$ws.Range("E23").Value = "  +9.39%  "
$ws.Range("D24").Value = "258.90"
$ws.Range("E24").Value = "  +13.08%  "
$ws.Range("D25").Value = "9.28"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "12.11"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "39.09"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").Value = "22.68"
$ws.Range("E30").Value = "  +7.92%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "173.50"
$ws.Range("B32").Value = "WEMIXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "0.0926"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("E34").Value = "  +5.41%  "
$ws.Range("D35").Value = "0.132"
$ws.Range("E35").Value = "  +5.18%  "
$ws.Range("D36").Value = "4.96"
$ws.Range("E36").Value = "  -4.59%  "
$ws.Range("E37").Value = "  -4.49%  "
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("E40").Value = "  +15.94%  "
$ws.Range("E41").Value = "  +14.05%  "
$ws.Range("D42").Value = "71.94"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").Value = "12.79"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("D47").Value = "9.37"
$ws.Range("E47").Value = "  +10.75%  "
$ws.Range("D48").Value = "111.64"
$ws.Range("E48").Value = "  +7.84%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("D51").Value = "0.472"
$ws.Range("E51").Value = "  +7.28%  "
